$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.302.96"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "3.511.91"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'610.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("D6").Value = "'150.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("D7").Value = "3.510.27"
$ws.Range("E7").Value = "  -1.08%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("D13").Value = "'0.0000220"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.27%  "

$ws.Range("D14").Value = "4.103.67"
$ws.Range("E14").Value = "  -1.25%  "

$ws.Range("D15").Value = "'31.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").Value = "3.516.35"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").Value = "67.293.35"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "'6.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "'15.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").Value = "'442.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").Value = "'9.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.61%  "

$ws.Range("D23").Value = "'0.625"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.19%  "

$ws.Range("D24").Value = "'77.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").Value = "'0.0000129"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.78%  "

$ws.Range("D26").Value = "3.650.39"
$ws.Range("E26").Value = "  -1.31%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'10.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.52%  "

$ws.Range("D29").Value = "'8.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("D31").Value = "'1.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.54%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("E33").Value = "  +2.99%  "

$ws.Range("D34").Value = "'25.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("D35").Value = "'6.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").Value = "3.501.44"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("E37").Value = "  -3.82%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'177.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "'2.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").Value = "'0.0872"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "

$ws.Range("D44").Value = "'5.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.36%  "

$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D46").Value = "'45.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").Value = "'27.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.31%  "

$ws.Range("E48").Value = "  +5.10%  "

$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("E50").Value = "  -1.78%  "

$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.55%  "
